$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14142
$ws1.Range("F3").Value = 332
$ws1.Range("F4").Value = 684
$ws1.Range("F6").Value = 539
$ws1.Range("F7").Value = 1478

# Sheet "全部类型" - update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14142
$ws4.Range("F3").Value = 332
$ws4.Range("F4").Value = 684
$ws4.Range("F8").Value = 539
$ws4.Range("F9").Value = 1478
